$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2801.8333
$ws.Range("I116").Value = 2851.25
$ws.Range("J116").Value = 2703
$ws.Range("K116").Value = 2851.25
$ws.Range("L116").Value = 2703
$ws.Range("M116").Value = 590.75
$ws.Range("N116").Value = -9587

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 69506.5
$ws.Range("J15").Value = 69506.5
$ws.Range("L15").Value = 69506.5
$ws.Range("N15").Value = -70206.5
$ws.Range("H22").Value = 2500
$ws.Range("J22").Value = 2500
$ws.Range("L22").Value = 2500
$ws.Range("N22").Value = -3098
$ws.Range("H133").Value = 26000
$ws.Range("J133").Value = 26000
$ws.Range("L133").Value = 26000
$ws.Range("N133").Value = -31060

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 7004
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H64").Value = 750
$ws.Range("J64").Value = 716.6667
$ws.Range("L64").Value = 716.6667
$ws.Range("N64").Value = -1166.6667
$ws.Range("H67").Value = 750
$ws.Range("J67").Value = 716.6667
$ws.Range("L67").Value = 716.6667
$ws.Range("N67").Value = -2276.6667
$ws.Range("H86").Value = 2740.8462
$ws.Range("I86").Value = 2643.25
$ws.Range("J86").Value = 2897
$ws.Range("K86").Value = 2643.25
$ws.Range("L86").Value = 2897
$ws.Range("M86").Value = -1520.25
$ws.Range("N86").Value = -5143
$ws.Range("H89").Value = 2740.8462
$ws.Range("I89").Value = 2643.25
$ws.Range("J89").Value = 2897
$ws.Range("K89").Value = 13216.25
$ws.Range("L89").Value = 14485
$ws.Range("M89").Value = -7600.25
$ws.Range("N89").Value = -25717
$ws.Range("H99").Value = 1919.3
$ws.Range("I99").Value = 1854.7778
$ws.Range("K99").Value = 1854.7778
$ws.Range("M99").Value = -356.7778000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 696.65
$ws.Range("I5").Value = 501.83334
$ws.Range("J5").Value = 2450
$ws.Range("K5").Value = 1505.50002
$ws.Range("L5").Value = 7350
$ws.Range("M5").Value = -1393.50002
$ws.Range("N5").Value = -7574
$ws.Range("H12").Value = 109.55556
$ws.Range("J12").Value = 125.35714
$ws.Range("L12").Value = 376.07142
$ws.Range("N12").Value = -722.07142
$ws.Range("H92").Value = 600
$ws.Range("I92").Value = 600
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1800
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -552
$ws.Range("N92").ClearContents()
$ws.Range("H98").Value = 339.25
$ws.Range("I98").Value = 203
$ws.Range("J98").Value = 384.66666
$ws.Range("K98").Value = 609
$ws.Range("L98").Value = 1153.99998
$ws.Range("M98").Value = 889
$ws.Range("N98").Value = -4149.999980000001
$ws.Range("H108").Value = 1675.6666
$ws.Range("I108").Value = 1013.5
$ws.Range("J108").Value = 3000
$ws.Range("K108").Value = 3040.5
$ws.Range("L108").Value = 9000
$ws.Range("M108").Value = -160.5
$ws.Range("N108").Value = -14760
$ws.Range("H109").Value = 3595.3333
$ws.Range("J109").Value = 7755
$ws.Range("L109").Value = 23265
$ws.Range("N109").Value = -25345
$ws.Range("H112").Value = 3658.8235
$ws.Range("J112").Value = 3153.3333
$ws.Range("L112").Value = 9459.999899999999
$ws.Range("N112").Value = -11675.9999
$ws.Range("H117").Value = 459.125
$ws.Range("J117").Value = 660.75
$ws.Range("L117").Value = 1982.25
$ws.Range("N117").Value = -8866.25
$ws.Range("H118").Value = 2974.4443
$ws.Range("J118").Value = 3034.524
$ws.Range("L118").Value = 9103.572
$ws.Range("N118").Value = -11589.572
$ws.Range("H121").Value = 1038.7
$ws.Range("I121").Value = 225.07692
$ws.Range("J121").Value = 1324.5676
$ws.Range("K121").Value = 675.23076
$ws.Range("L121").Value = 3973.7028
$ws.Range("M121").Value = 634.76924
$ws.Range("N121").Value = -6593.7028
$ws.Range("H122").Value = 6715
$ws.Range("I122").Value = 470.9
$ws.Range("J122").Value = 15635.143
$ws.Range("K122").Value = 4238.099999999999
$ws.Range("L122").Value = 140716.287
$ws.Range("M122").Value = -1788.099999999999
$ws.Range("N122").Value = -145616.287
$ws.Range("H123").Value = 6781
$ws.Range("I123").Value = 3030
$ws.Range("J123").Value = 7249.875
$ws.Range("K123").Value = 9090
$ws.Range("L123").Value = 21749.625
$ws.Range("M123").Value = -6640
$ws.Range("N123").Value = -26649.625
$ws.Range("H124").Value = 101828.6
$ws.Range("I124").Value = 334500
$ws.Range("J124").Value = 2112.2856
$ws.Range("K124").Value = 1003500
$ws.Range("L124").Value = 6336.8568
$ws.Range("M124").Value = -998590
$ws.Range("N124").Value = -16156.8568
$ws.Range("H131").Value = 929.5599999999999
$ws.Range("I131").Value = 576.2727
$ws.Range("J131").Value = 1207.1428
$ws.Range("K131").Value = 1728.8181
$ws.Range("L131").Value = 3621.4284
$ws.Range("M131").Value = 3311.1819
$ws.Range("N131").Value = -13701.4284
$ws.Range("H132").Value = 2734.9312
$ws.Range("I132").Value = 2792.3333
$ws.Range("J132").Value = 2719.9565
$ws.Range("K132").Value = 25130.9997
$ws.Range("L132").Value = 24479.6085
$ws.Range("M132").Value = -22600.9997
$ws.Range("N132").Value = -29539.6085
$ws.Range("H133").Value = 13967.385
$ws.Range("I133").Value = 7906
$ws.Range("J133").Value = 17755.75
$ws.Range("K133").Value = 23718
$ws.Range("L133").Value = 53267.25
$ws.Range("M133").Value = -18658
$ws.Range("N133").Value = -63387.25
$ws.Range("H134").Value = 7747.7646
$ws.Range("I134").Value = 5412
$ws.Range("J134").Value = 8721
$ws.Range("K134").Value = 16236
$ws.Range("L134").Value = 26163
$ws.Range("M134").Value = -11166
$ws.Range("N134").Value = -36303
$ws.Range("H135").Value = 696.65
$ws.Range("I135").Value = 501.83334
$ws.Range("J135").Value = 2450
$ws.Range("K135").Value = 4516.50006
$ws.Range("L135").Value = 22050
$ws.Range("M135").Value = -1981.50006
$ws.Range("N135").Value = -27120
$ws.Range("H136").Value = 1723.8
$ws.Range("I136").Value = 1154.75
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 3464.25
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = 1635.75
$ws.Range("N136").Value = -22200
$ws.Range("H137").Value = 10426899
$ws.Range("I137").Value = 27795022
$ws.Range("J137").Value = 6025.4
$ws.Range("K137").Value = 83385066
$ws.Range("L137").Value = 18076.2
$ws.Range("M137").Value = -83379966
$ws.Range("N137").Value = -28276.2

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 27779712
$ws.Range("I82").Value = 55557956
$ws.Range("J82").Value = 1466.6666
$ws.Range("K82").Value = 55557956
$ws.Range("L82").Value = 1466.6666
$ws.Range("M82").Value = -55557595
$ws.Range("N82").Value = -2188.6666
$ws.Range("H85").Value = 27779712
$ws.Range("I85").Value = 55557956
$ws.Range("J85").Value = 1466.6666
$ws.Range("K85").Value = 55557956
$ws.Range("L85").Value = 1466.6666
$ws.Range("M85").Value = -55556708
$ws.Range("N85").Value = -3962.6666
$ws.Range("H93").Value = 13823.25
$ws.Range("I93").Value = 13823.25
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 13823.25
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -12575.25
$ws.Range("N93").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H75").Value = 46600
$ws.Range("J75").Value = 46600
$ws.Range("L75").Value = 46600
$ws.Range("N75").Value = -48472
$ws.Range("H78").Value = 46600
$ws.Range("J78").Value = 46600
$ws.Range("L78").Value = 139800
$ws.Range("N78").Value = -149160
$ws.Range("H132").Value = 4763559.5
$ws.Range("I132").Value = 1579.7894
$ws.Range("J132").Value = 10418411
$ws.Range("K132").Value = 4739.3682
$ws.Range("L132").Value = 31255233
$ws.Range("M132").Value = -2209.3682
$ws.Range("N132").Value = -31260293
